{"js": "// Use case adjustment: correct a team member's name in the \"ALUNOS\" (students)\n// list of the requirements document: \"Pedro Henrique\" -> \"Tailyne Bertoncelli\".\n\nconst body = context.document.body;\n\n// --- 1. Replace the student name -------------------------------------------------\nconst nameResults = body.search(\"Pedro Henrique\", { matchCase: true });\nnameResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n    nameResults.items[0].insertText(\"Tailyne Bertoncelli\", Word.InsertLocation.replace);\n    await context.sync();\n}\n\n// --- 2. Drop a collapsed \"_GoBack\" bookmark right after the edit location --------\n// (mirrors Word's own behaviour of remembering the last edit position on save)\nconst anchorResults = body.search(\"Cherman,\", { matchCase: true });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n    const existing = body.getBookmarkRangeOrNullObject(\"_GoBack\");\n    existing.load(\"isNullObject\");\n    await context.sync();\n    if (!existing.isNullObject) {\n        context.document.deleteBookmark(\"_GoBack\");\n        await context.sync();\n    }\n\n    const goBackPoint = anchorResults.items[0].getRange(Word.RangeLocation.end);\n    goBackPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n}\n", "ps1": "# Use case adjustment: correct a team member's name in the \"ALUNOS\" (students)\n# list of the requirements document: \"Pedro Henrique\" -> \"Tailyne Bertoncelli\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Replace the student name -------------------------------------------------\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"Pedro Henrique\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    $findRange.Text = \"Tailyne Bertoncelli\"\n}\n\n# --- 2. Drop a collapsed \"_GoBack\" bookmark right after the edit location --------\n# (mirrors Word's own behaviour of remembering the last edit position on save)\n$markRange = $d.Content\n$markFind = $markRange.Find\n$markFind.Text = \"Cherman,\"\n$markFind.MatchCase = $true\n$markFound = $markFind.Execute()\n\nif ($markFound) {\n    $goBackRange = $d.Range($markRange.End, $markRange.End)\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n    $d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n}\n"}
